$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sheet1" to "Results"
$ws.Name = "Results"

# Update the "Reason" note for the training run in row 17 to reflect that
# this configuration is being kept due to its lower loss.
$ws.Range("Y17").Value = "Pretty much as the previous good one. Keeping this one due to lower loss."

# Move the active selection to A20 (was D23)
$ws.Range("A20").Select()
